$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "0.353 (0.056)"
$ws.Range("D2").Value = "0.330 (0.037)"
$ws.Range("E2").Value = "0.184 (0.042)"

$ws.Range("C3").Value = "0.390 (0.009)"
$ws.Range("D3").Value = "0.346 (0.009)"
$ws.Range("E3").Value = "0.194 (0.001)"

$ws.Range("C4").Value = "0.453 (0.020)"
$ws.Range("D4").Value = "0.404 (0.019)"
$ws.Range("E4").Value = "0.197 (0.006)"

$ws.Range("C5").Value = "0.467 (0.018)"
$ws.Range("D5").Value = "0.418 (0.019)"
$ws.Range("E5").Value = "0.204 (0.004)"

$ws.Range("C6").Value = "0.124 (0.085)"
$ws.Range("D6").Value = "0.137 (0.101)"

$ws.Range("C7").Value = "0.128 (0.081)"
$ws.Range("D7").Value = "0.150 (0.087)"
$ws.Range("E7").Value = "0.014 (0.010)"

$ws.Range("C8").Value = "0.208 (0.001)"
$ws.Range("D8").Value = "0.238 (0.002)"
$ws.Range("E8").Value = "0.006 (0.007)"

$ws.Range("C9").Value = "0.300 (0.049)"
$ws.Range("D9").Value = "0.301 (0.045)"
$ws.Range("E9").Value = "0.106 (0.025)"
